$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = 'Parneet Kaur'

# Preconditions (column E), filled top-to-bottom
$ws.Range("E7").Value = 'None'
$ws.Range("E8").Value = 'None'
$ws.Range("E9").Value = 'None'
$ws.Range("E10").Value = 'None'
$ws.Range("E11").Value = 'Valid BankAccount object exists.'
$ws.Range("E12").Value = 'Valid BankAccount object exists.'
$ws.Range("E13").Value = 'Valid BankAccount object exists.'
$ws.Range("E14").Value = 'Start balance = 100.00.'
$ws.Range("E15").Value = 'Start balance = 100.00.'
$ws.Range("E16").Value = 'Start balance = 50.00.'
$ws.Range("E17").Value = 'Start balance = 100.00.'
$ws.Range("E18").Value = 'None'
$ws.Range("E19").Value = 'Start balance = 200.00.'
$ws.Range("E20").Value = 'None'
$ws.Range("E21").Value = 'Start balance = 50.00.'
$ws.Range("E22").Value = 'Valid BankAccount object exists.'

# Method Inputs (column F), filled top-to-bottom
$ws.Range("F7").Value = '(20019, 1010, 100.00)'
$ws.Range("F8").Value = '(20019, 1010, "abc")'
$ws.Range("F9").Value = '("X", 1010, 0.0)'
$ws.Range("F10").Value = '(20019, "X", 0.0)'
$ws.Range("F11").Value = 'obj.account_number'
$ws.Range("F12").Value = 'obj.client_number'
$ws.Range("F13").Value = 'obj.balance'
$ws.Range("F14").Value = 'obj.update_balance(25.50)'
$ws.Range("F15").Value = 'obj.update_balance(-40)'
$ws.Range("F16").Value = 'obj.update_balance("abc")'
$ws.Range("F17").Value = 'obj.deposit(75.34)'
$ws.Range("F18").Value = 'obj.deposit(0)'
$ws.Range("F19").Value = 'obj.withdraw(75.34)'
$ws.Range("F20").Value = 'obj.withdraw(0)'
$ws.Range("F21").Value = 'obj.withdraw(75)'
$ws.Range("F22").Value = 'str(obj)'

# Expected Result (column G), filled top-to-bottom
$ws.Range("G7").Value = 'Object created successfully. All private attributes set to correct values.'
$ws.Range("G8").Value = 'Balance initialized to 0.00.'
$ws.Range("G9").Value = 'ValueError raised: “Account number must be an integer.”'
$ws.Range("G10").Value = 'ValueError raised: “Client number must be an integer.”'
$ws.Range("G11").Value = 'Returns integer account number (e.g. 20019).'
$ws.Range("G12").Value = 'Returns integer client number (e.g. 1010).'
$ws.Range("G13").Value = 'Returns float balance (e.g. 100.00).'
$ws.Range("G14").Value = 'Balance updated to 125.50.'
$ws.Range("G15").Value = 'Balance updated to 60.00.'
$ws.Range("G16").Value = 'Balance remains 50.00.'
$ws.Range("G17").Value = 'Balance updated to 175.34.'
$ws.Range("G18").Value = 'ValueError raised: “Deposit amount: $0.00 must be positive.”'
$ws.Range("G19").Value = 'Balance updated to 124.66.'
$ws.Range("G20").Value = 'ValueError raised: “Withdraw amount: $0.00 must be positive.”'
$ws.Range("G21").Value = 'ValueError raised: “Withdraw amount: $75.00 must not exceed the account balance: $50.00.”'
$ws.Range("G22").Value = 'Returns: "Account Number: 20019 Balance: $6,764.67" followed by newline.'

# View state: zoom + selection (matches final author session)
$excel.ActiveWindow.Zoom = 67
$ws.Range("K21").Select()
